# EMC: - dalsi slajdy
# Add a third data column (D) to the "Antena" sheet: D2 = 50,
# D3 = 1e-6 (scientific-formatted, like column B), and
# D4 = D3/(D2^2) (also scientific-formatted).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Antena")

$ws.Range("D2").Value = 50

$ws.Range("D3").Value = 0.000001
$ws.Range("D3").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("D4").Formula = "=D3/(D2^2)"
$ws.Range("D4").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Activate()
$ws.Range("D4").Select()
